$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the shared string used by D2: "Bisual" -> "Visual"
$ws.Range("D2").Value = "Attribute/Sensory/Visual"

# Replace the long HED tag text in D3 with the shorter "Awake"
$ws.Range("D3").Value = "Awake"

# Drop the ResponseOnset / ResponseOffset rows (4 and 5) - clear their contents
$ws.Range("A4:D5").ClearContents()

# Row heights shrink now that the long text is gone
$ws.Rows.Item(3).RowHeight = 14.9
$ws.Rows.Item(5).RowHeight = 13.8

# Move selection/scroll back to the top-left of the sheet
[void]$ws.Range("A1").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
[void]$ws.Range("A4").Select()
